$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices and 1h volume changes),
# including a newly-inserted "RocketPoolETH" row that shifts rows 44-51 down by one.

$ws.Range("D2").Value = "29.397.83"
$ws.Range("D3").Value = "1.848.68"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.29"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6298"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07622"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2941"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07746"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.842.69"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.009"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001086"
$ws.Range("E14").Value = "  +8.31%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.46"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "2.093.50"
$ws.Range("E17").Value = "  -7.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.136"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "29.430.65"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.50"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.445"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.22"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.377"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.470"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.300"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05630"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.114"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.043"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.853"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.586"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.779"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Value = "1.229.81"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01797"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.479"
$ws.Range("E41").Value = "  +4.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9089"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.002.75"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.41"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.11"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("E47").Value = "  +3.63%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.156"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4012"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.062"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.686"
$ws.Range("E51").Value = "  -0.51%  "
